# Session Info workbook update — "updated some session data and debugged webhook"
#
# Content changes (row numbers refer to the worksheet, not the Session ID
# column, which is one less than the row number):
#   H1  : " Location" (leading space)                 -> "Location"
#   G13 : "Discussing how AI capabilities..."          -> "Discussing how A I capabilities..."
#   G14 : "...government CXOs have become..."          -> "...government C X Os have become..."
#   G26 : "What if ROI was ... Central IT ... IT operations ... Public ROI? ..."
#         -> "What if R O I was ... Central I T ... I T operations ... Public R O I? ..."

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Session 12 ("AI in Mission Decision Making") description: spell out "AI".
$ws.Range("G13").Value = "Discussing how A I capabilities can provide effective distillation of information at scale for improved and timely decision making."

# Session 13 ("Driving Innovation to Support the Mission") description: spell out "CXOs".
$ws.Range("G14").Value = "Innovation comes in many forms, including broth both process-orient and technical.  with budgets increasinly tight and an ever-expanding workload, government C X Os have become the drivers of innovation across the Federal enterprise.  Join us for adiscussion on how a diverse group of leaders from across government are working to foster a culture that promotes and rewards innovation, drives change, and benefits mission performance by empowering employees across their respective organizations."

# Header: "Location" column label loses its stray leading space.
$ws.Range("H1").Value = "Location"

# Session 25 ("Why the Future of Government is a Citizen Experience") description: spell out "ROI"/"IT".
$ws.Range("G26").Value = "What if R O I was not in a spreadsheet or financial ledger, but it was the constituent experience and the satisfaction of a job done well. Not just well done even - amazingly done. We, with a partnership of industry and public sector, must focus on the experience of the actual people we need to serve and reduce the burden on our public sector workforce.  State governments are harnessing the momentum of change to focus on the citizen experience and enabling agencies to innovate and modernize the way we serve customers. State agencies are propelling the path to production, while Central I T is utilizing strategic technologies to empower agencies to act faster and mitigate the unforeseen risk inherent in I T operations.  How do we realize Public R O I? By removing the barriers to innovation, owning the hurdles that slow operations down and embracing the momentum of change."

# Leave the selection where the author's last edit was: G26 (also scrolled the
# viewport down a bit while working through rows 21-26).
$ws.Range("G26").Select()
